# QA ONB test and L1 regression
#
# This script reproduces (via Excel COM automation) the edit that:
#   - adds a new data row (24th iteration, r=29) to "Foglio1", carrying the
#     corrected "Matrix non positive def" label (fixing the old typo'd
#     "Matrix non positve def" references in the same column on row 20),
#   - leaves "Foglio2" values untouched, and
#   - ends the session with "Foglio1" as the active/selected sheet
#     (instead of "Foglio2"), scrolled near the new row, with H30 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws2 = $wb.Worksheets.Item("Foglio2")

# --- Foglio1: fix the mislabeled "Matrix non positve def" entries on row 20
#     so they read the corrected "Matrix non positive def" text -----------
$ws1.Range("E20").Value = "Matrix non positive def"
$ws1.Range("F20").Value = "Matrix non positive def"

# --- Foglio1: append the new iteration row (r=29) ------------------------
# Copy the formatting of the row above (so A29 keeps the same "Agency FB"
# centered numbering style, and the row picks up the usual 15pt row
# height) before filling in the new values.
$ws1.Range("A28:I28").Copy()
$ws1.Range("A29:I29").PasteSpecial(-4122)
$ws1.Rows.Item(29).RowHeight = 15

$ws1.Range("A29").Value = 24
$ws1.Range("B29").Value = 13
$ws1.Range("C29").Value = 6
$ws1.Range("D29").Value = 6
$ws1.Range("E29").Value = "Matrix non positive def"
$ws1.Range("F29").Value = "Matrix non positive def"
$ws1.Range("G29").Value = "Matrix non positive def"
$ws1.Range("H29").Value = 8
$ws1.Range("I29").Value = 11

# --- Switch the active sheet/selection from Foglio2 back to Foglio1 ------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("H30").Select()
